$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jobDescription = @'
Please find the Job Description (JD) below for your reference.
Job Role: ML Engineer/Data Scientist
Job Description:
• Develop, train, test, and deploy machine learning models across various fields including computer vision, LLMs, and with tabular and time series data.
• Strong experience in Python, FastAPI, Flask
• Strong experience in SQL
• Strong experience Design Pattern/algorithms and data structures
• Familiarity with OOPS, Design Pattern/algorithms and data structures
• Familiarity with continuous integration, deployment, and automated build processes for scalable application delivery using Docker/Kubernetes
• Practical knowledge of one or more major cloud platforms (e.g. Azure, AWS, or GCP).
• Excellent written and verbal communication skills in English.
• Experiment with novel deep learning-based technologies such as self-supervised learning and generative AI. 
• Work directly with customer data and set up data pipelines to collect, curate, transform, and version data. 
• Participate in the collection, analysis, interpretation, and output of large amounts of data using advanced AI techniques like deep learning, NLP, and computer vision good foundational experience in PyTorch / Tensorflow.
• Work within the global corporate Artificial Intelligence division, which addresses real business challenges and opportunities across multiple countries.
• Collaborate across different business and corporate functions in an international team composed of Project Managers, Data Scientists, Data and Software Engineers within the Artificial Intelligence team and others in the Global AI team
 Qualifications:
• Bachelor’s degree or master’s degree in data science, Computational Statistics/Mathematics, Computer Science, or related field
• Fluent English
'@

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Fullstack Developer"
$ws.Range("C3").Value = $jobDescription
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4
